# Update the "Turkey 1 Lig" sheet: several rows had their match data
# (everything except the running index in column A) shuffled between rows.
# For every group of affected rows, the content of columns B:AD moves to a
# different row within the same group (pairwise swaps, plus a couple of
# 3-row rotations). Column A (the sequential id) always stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Groups of row numbers whose B:AD content gets rearranged among themselves.
$groups = @(
  @(22,23),
  @(85,86),
  @(97,98),
  @(103,104),
  @(122,123),
  @(141,142),
  @(169,170),
  @(205,206),
  @(220,221),
  @(241,242),
  @(246,247),
  @(295,296,297),
  @(300,305,306)
)

# For a given row r, $rowMap[r] is the row whose ORIGINAL B:AD content
# should end up in row r after the edit.
$rowMap = @{
  22=23; 23=22;
  85=86; 86=85;
  97=98; 98=97;
  103=104; 104=103;
  122=123; 123=122;
  141=142; 142=141;
  169=170; 170=169;
  205=206; 206=205;
  220=221; 221=220;
  241=242; 242=241;
  246=247; 247=246;
  295=297; 296=295; 297=296;
  300=305; 305=306; 306=300
}

foreach ($grp in $groups) {
  # Snapshot the original B:AD values of every row in the group first,
  # so rotations (3+ rows) don't clobber data before it is read.
  $orig = @{}
  foreach ($r in $grp) {
    $orig[$r] = $ws.Range("B$r`:AD$r").Value2
  }

  foreach ($r in $grp) {
    $srcRow = $rowMap[$r]
    $ws.Range("B$r`:AD$r").Value2 = $orig[$srcRow]
  }
}
